$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = -22.06139999999999
$ws.Range("A12").Value = -21.39669999999999
$ws.Range("A18").Value = -22.11000000000001
$ws.Range("A37").Value = -19.77829999999999
$ws.Range("A55").Value = -22.1976
$ws.Range("A68").Value = -21.44299999999999
$ws.Range("A77").Value = -20.45179999999999
$ws.Range("A78").Value = -19.72659999999998
